# Apply the "add various YFCF files" edit to Sheet1 of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Widen column A (stores as width="14.5" in the OOXML <col> element; Excel's
# COM ColumnWidth (characters) is offset from the stored width by 5/6).
$ws.Columns("A").ColumnWidth = 13.666666666666666

# New row 30: a lone label in column A.
$ws.Range("A30").Value = "\YFCF"

# New row 32: SonicYFCF entry.
$ws.Range("A32").Value = "SonicYFCF"
$ws.Range("B32").Value = "https://web.archive.org/web/19991104065234im_/http://dewey.rug.ac.be/YFCF/SonicYFCF.html"

# New row 33: HotSonic entry.
$ws.Range("A33").Value = "HotSonic"
$ws.Range("B33").Value = "https://web.archive.org/web/19991104070439im_/http://dewey.rug.ac.be/YFCF/HotSonic.html"

# New row 34: UnploughCD entry - the URL was filled in before the label.
$ws.Range("B34").Value = "https://web.archive.org/web/19991104070723im_/http://dewey.rug.ac.be/YFCF/UnploughCD.html"
$ws.Range("A34").Value = "UnploughCD"

# Leave the selection where the author's cursor ended up.
$ws.Range("B36").Select()
